$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows of grind stats to append starting at row 94
$newRows = @(
    @("2026-02-26 17:59:55", 86344, 102063, 27752, 50654, 31079, 17393, 0, 19136, 55513, 493, 4438, 0, 0, 151, 172, 1800000, 635335, 0),
    @("2026-02-26 18:06:26", 89575, 106260, 29014, 52186, 32087, 17844, 0, 19752, 57935, 493, 4633, 0, 0, 152, 173, 0, 718754, 0),
    @("2026-02-26 19:06:06", 93343, 111715, 29834, 53964, 33585, 18109, 810, 20639, 60711, 493, 4870, 0, 0, 154, 175, 1800000, 812874, 0),
    @("2026-02-26 19:13:12", 97420, 115299, 31307, 56694, 35051, 19696, 812, 21489, 62871, 493, 5135, 0, 0, 154, 176, 1800000, 914296, 0)
)

$startRow = 94
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowIndex = $startRow + $i
    $rowData = $newRows[$i]

    # Column A holds the timestamp as plain text
    $ws.Cells.Item($rowIndex, 1).Value = $rowData[0]

    for ($col = 1; $col -lt $rowData.Count; $col++) {
        $ws.Cells.Item($rowIndex, $col + 1).Value = $rowData[$col]
    }
}
